# Update countries & provincias Spain
#
# 1) Barein overtakes Kazajistan in the ranking, so the two country rows
#    (59 and 60) swap places; Barein's row gets freshly updated stats.
# 2) Refresh the day's case numbers for Estados Unidos (row 4) and
#    Ecuador (row 25).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the Kazajistan / Barein country labels (row 59 <-> row 60) ---
$ws.Range("A59").Value = "Barein"
$ws.Range("A60").Value = "Kazajistan"

# --- Row 59 (now Barein) gets the updated counts ---
$ws.Range("B59").Value = 5236
$ws.Range("C59").Value = 295
$ws.Range("D59").Value = 2152
$ws.Range("E59").Value = 3076
$ws.Range("F59").Value = 6
$ws.Range("G59").Value = 0
$ws.Range("H59").Value = 8

# --- Row 60 (now Kazajistan) takes over what used to be row 59's counts ---
$ws.Range("B60").Value = 5207
$ws.Range("C60").Value = 117
$ws.Range("D60").Value = 2074
$ws.Range("E60").Value = 3101
$ws.Range("F60").Value = 33
$ws.Range("G60").Value = 1
$ws.Range("H60").Value = 32

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 1380443
$ws.Range("C4").Value = 12805
$ws.Range("E4").Value = 1039915
$ws.Range("F4").Value = 16465
$ws.Range("G4").Value = 668
$ws.Range("H4").Value = 81455

# --- Row 25: Ecuador ---
$ws.Range("B25").Value = 29509
$ws.Range("E25").Value = 23931
$ws.Range("G25").Value = 18
$ws.Range("H25").Value = 2145
